$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 123, shifting existing rows 123..176 down to 124..177.
$ws.Rows.Item(123).Insert()

# Populate the newly inserted row 123 with the new weekly data point.
$ws.Cells.Item(123, 1).Value = 11
$ws.Cells.Item(123, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(123, 3).Value = "Bíobío"
$ws.Cells.Item(123, 4).Value = 44726
$ws.Cells.Item(123, 5).Value = 8
$ws.Cells.Item(123, 6).Value = "Fruta"
$ws.Cells.Item(123, 7).Value = 100108
$ws.Cells.Item(123, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(123, 9).Value = 100108005
$ws.Cells.Item(123, 10).Value = "Piña"
$ws.Cells.Item(123, 11).Value = "Caramelo"
$ws.Cells.Item(123, 12).Value = "Segunda"
$ws.Cells.Item(123, 13).Value = 200
$ws.Cells.Item(123, 14).Value = 17000
$ws.Cells.Item(123, 15).Value = 18000
$ws.Cells.Item(123, 16).Value = 17500
$ws.Cells.Item(123, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(123, 18).Value = "Ecuador"
$ws.Cells.Item(123, 19).Value = 1250
$ws.Cells.Item(123, 20).Value = 14
